$d = $word.ActiveDocument

# --- Change 1: merge the two "SAT Oct 27" / " 11:25:55 IST 2018" runs
#     into a single run with the combined text. A find/replace with
#     identical replacement text collapses the run split.
$d.Content.Find.Execute("SAT Oct 27 11:25:55 IST 2018", $true, $false, $false, $false, $false, $true, 1, $false, "SAT Oct 27 11:25:55 IST 2018", 2) | Out-Null

# --- Change 2: append the new "SUN Oct 28" chick-in block right after
#     the last "Amount Received mode ... - CASH" paragraph, before the
#     pre-existing blank paragraphs.

# Locate that paragraph (walk from the end, it's the last one containing "- CASH").
$paras = $d.Paragraphs
$anchor = $null
for ($i = $paras.Count; $i -ge 1; $i--) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*- CASH*") {
        $anchor = $p
        break
    }
}

# 1) blank paragraph
$anchor.Range.InsertParagraphAfter()
$p1 = $anchor.Next()

# 2) "SUN Oct 28 12:37:49 IST 2018"
$p1.Range.InsertParagraphAfter()
$p2 = $p1.Next()
$p2.Range.InsertAfter("SUN Oct 28")
$p2.Range.InsertAfter(" 12:37:49 IST 2018")

# 3) "Person Name" ... "- DHEVEGOWDA"
$p2.Range.InsertParagraphAfter()
$p3 = $p2.Next()
$p3.Range.InsertAfter("Person Name`t`t`t`t- DHEVEGOWDA")

# 4) dashed separator line
$p3.Range.InsertParagraphAfter()
$p4 = $p3.Next()
$p4.Range.InsertAfter("---------------------------------------------------------------")

# 5) "Item Name" ... "- BEET"
$p4.Range.InsertParagraphAfter()
$p5 = $p4.Next()
$p5.Range.InsertAfter("Item Name`t`t`t`t- BEET")

# 6) "Amount Received" ... "- 2150" (red)
$p5.Range.InsertParagraphAfter()
$p6 = $p5.Next()
$p6.Range.InsertAfter("Amount Received`t`t`t- 2150")
$p6.Range.Font.Color = 255

# 7) "Amount Received mode" ... "- CASH AND CLEARD"
$p6.Range.InsertParagraphAfter()
$p7 = $p6.Next()
$p7.Range.InsertAfter("Amount Received mode`t`t- CASH AND CLEARD")

# 8) blank paragraph
$p7.Range.InsertParagraphAfter()
$p8 = $p7.Next()

# 9) blank paragraph
$p8.Range.InsertParagraphAfter()
